$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): MinRollCrn / MaxRollCrn
$ws.Range("D1").Value = "MinRollCrn"
$ws.Range("E1").Value = "MaxRollCrn"

# New data for columns D (MinRollCrn) and E (MaxRollCrn), rows 2-8
$ws.Range("D2").Value = -0.9
$ws.Range("E2").Value = 0.3

$ws.Range("D3").Value = -0.9
$ws.Range("E3").Value = 0.3

$ws.Range("D4").Value = -0.9
$ws.Range("E4").Value = 0.3

$ws.Range("D5").Value = -0.9
$ws.Range("E5").Value = 0.3

$ws.Range("D6").Value = -0.5
$ws.Range("E6").Value = 0.3

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

# Widen new column D to fit the "MinRollCrn" header (closest achievable snap to ~11.22 chars)
$ws.Columns.Item(4).ColumnWidth = 10.5

# Move the active selection down to the last data row/col, like in the edited file
$ws.Range("E8").Select()

# Page setup: printed on A4-ish "Letter"/paper id 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
